$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")
$summary = $wb.Worksheets.Item("Daily Summary")

# Insert a new row at row 2, pushing existing rows down
$ws.Rows.Item(2).Insert()

# Fill in the new order row
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "2026-01-13 18:56"
$ws.Cells.Item(2, 3).Value = "Sagar Borse"
$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "7588930329"
$ws.Cells.Item(2, 5).Value = "Test3,"
$ws.Cells.Item(2, 6).Value = "Square Heat Pad x1"
$ws.Cells.Item(2, 7).Value = 50
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"

# Update Daily Summary sheet
$summary.Cells.Item(2, 2).Value = 8
$summary.Cells.Item(2, 5).Value = 375
$summary.Cells.Item(2, 7).Value = 375
